# Updated all workflow files with new path to NVIS Extant and Pre1750 files
#
# The "RawDataPath" cell for the Terrestrial-Extant-IUCNGET row (row 4)
# pointed at an intermediate NVIS/IUCN-GET raster stamped 20240709; the
# upstream processing step now produces a file stamped 20240801, so the
# workbook's recorded path needs to be refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 = Terrestrial-Extant-IUCNGET; column B = RawDataPath.
$ws.Range("B4").Value = "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\processing\NEAP_intermediate\NVIS_IUCNGET_DK_20240801.tif"

# Leave the active selection on B5, matching the workbook's last saved UI
# state.
$ws.Range("B5").Select() | Out-Null
